# Insert a new data row at row 358 (shifting the existing rows 358-459 down
# to 359-460) and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 358, pushing everything below it down by one.
$ws.Rows.Item(358).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A358").Value = 6
$ws.Range("B358").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C358").Value = "Metropolitana"
$ws.Range("D358").Value = 44855
$ws.Range("E358").Value = 13
$ws.Range("F358").Value = 100112032
$ws.Range("G358").Value = "Zapallo italiano"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 800
$ws.Range("K358").Value = 12000
$ws.Range("L358").Value = 13000
$ws.Range("M358").Value = 12438
$ws.Range("N358").Value = "`$/caja 50 unidades"
$ws.Range("O358").Value = "Región de O'Higgins"
$ws.Range("P358").Value = 249
$ws.Range("Q358").Value = 50
$ws.Range("R358").Value = "Hortaliza"
